# Add spacing / formatting between each fio benchmark script section:
# one row per burst/sustain test with its name in column A, and a header
# row (row 5) above them describing the measured columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per fio test scenario, column A, rows 6-13.
$tests = @(
    "BURST SEQUENTIAL WRITE",
    "BURST SEQUENTIAL READ",
    "BURST RANDOM WRITE",
    "BURST RANDOM READ",
    "BURST RANDOM WRITE OIO",
    "BURST RANDOM READ OIO",
    "SUBSTAIN SEQUENTIAL WRITE",
    "SUBSTAIN SEQUENTIAL WRITE"
)

for ($i = 0; $i -lt $tests.Length; $i++) {
    $row = 6 + $i
    $ws.Cells.Item($row, 1).Value = $tests[$i]
}

# Column headers for the benchmark metrics, starting at column B of row 5.
$headers = @(
    "BLOCK SIZE",
    "IO-DEPTHS",
    "THREADS",
    "SIZE",
    "PRECONDITION",
    "IOPS",
    "BANDWIDTH",
    "AVG, LATENCY",
    "50th",
    "75th",
    "99th",
    "99.9th",
    "99.99th",
    "99.999th",
    "99.9999th"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # column B == 2
    $ws.Cells.Item(5, $col).Value = $headers[$i]
}
